$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C16").Value = "a"
$ws.Hyperlinks.Add($ws.Range("C16"), "https://example.com")
# try clearing alignment override
$ws.Range("C16").HorizontalAlignment = -4131
